# Apply the ValueSet/CodeSystem URL, version, date, and publisher updates
# (ibm.com -> linuxforhealth.org, 7.0.0 -> 8.0.0, date bump, Alvearie Team -> LinuxForHealth Team)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# URL (B2)
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/eng-communication-request-status"

# Version (B3)
$meta.Range("B3").Value = "8.0.0"

# Date (B8)
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher (B9)
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet 2: "Include from Engagement Commu" ---
$codes = $wb.Worksheets.Item("Include from Engagement Commu")

# System URI (B4)
$codes.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/eng-communication-request-status-reason"
